$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D, J, K, L, M, P
# Derived from a weekly re-sort of the underlying date-ordered dataset
$rowData = @{
    2 = @(44292, 40, 3000, 3000, 3000, 1000)
    3 = @(44574, 50, 3000, 3000, 3000, 1000)
    4 = @(44193, 70, 3000, 3000, 3000, 1000)
    5 = @(44804, 85, 3000, 3000, 3000, 1000)
    6 = @(44340, 54, 3000, 3000, 3000, 1000)
    7 = @(44390, 50, 3000, 3000, 3000, 1000)
    8 = @(44536, 125, 2200, 2200, 2200, 733)
    9 = @(44756, 104, 2800, 3000, 2904, 968)
    10 = @(44291, 45, 3000, 3000, 3000, 1000)
    11 = @(44165, 68, 3000, 3000, 3000, 1000)
    12 = @(44223, 80, 2500, 3000, 2781, 927)
    13 = @(44559, 68, 2000, 2000, 2000, 667)
    14 = @(44557, 104, 2000, 2500, 2260, 753)
    15 = @(44389, 81, 2800, 3000, 2889, 963)
    16 = @(44537, 88, 2000, 2200, 2091, 697)
    17 = @(44242, 95, 2500, 3000, 2737, 912)
    18 = @(44187, 65, 3000, 3000, 3000, 1000)
    19 = @(44166, 45, 2500, 2500, 2500, 833)
    20 = @(44669, 92, 2500, 3000, 2755, 918)
    21 = @(44179, 78, 3000, 3000, 3000, 1000)
    22 = @(44221, 50, 2500, 2500, 2500, 833)
    23 = @(44225, 56, 3000, 3000, 3000, 1000)
    24 = @(44260, 60, 3500, 3500, 3500, 1167)
    25 = @(44222, 45, 3000, 3000, 3000, 1000)
    26 = @(44627, 78, 3500, 3500, 3500, 1167)
    27 = @(44224, 67, 3000, 3000, 3000, 1000)
    28 = @(44243, 45, 3000, 3000, 3000, 1000)
    29 = @(44845, 80, 2500, 2500, 2500, 833)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $vals[1]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals[2]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals[3]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals[5]   # P - Precio $/Kg
}
